$d = $word.ActiveDocument

# 1) "We have previously written a method " -> "Previously written method "
$d.Content.Find.Execute(
    "We have previously written a method ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Previously written method ", 2) | Out-Null

# 2) ", whose specification appears at the bottom of this page. It "
#    -> ", whose specification appears at the bottom of this page, "
$d.Content.Find.Execute(
    ", whose specification appears at the bottom of this page. It ", $true, $false, $false, $false, $false,
    $true, 1, $false, ", whose specification appears at the bottom of this page, ", 2) | Out-Null

# 3) "sorts an " -> "sorts " ; remember where the replaced text ends so the
#    _GoBack bookmark can be relocated there afterwards.
$sortsRange = $d.Content
$sortsRange.Find.Execute(
    "sorts an ", $true, $false, $false, $false, $false,
    $true, 1, $false, "sorts ", 2) | Out-Null
$sortsRange.Collapse(0)
$goBackStart = $sortsRange.Start
$goBackEnd = $sortsRange.End

# 4) Move the (hidden) _GoBack bookmark from the "radix" paragraph to right
#    after "sorts " in the paragraph above, matching the edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackRange = $d.Range($goBackStart, $goBackEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
